$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (2023-09-12 -> 2023-09-13, i.e. 45181 -> 45182) for every data row.
$ws.Range("C2:C257").Value = 45182
